$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 4: "Tool-basedCodeAnalysis"  -- SonarQube for IDE review + findings
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Tool-basedCodeAnalysis")

$ws4.Range("I3").Value = "Ticu Cristian"
$ws4.Range("J3").Value = 237
$ws4.Range("D4").Value = "SonarQube for IDE"

# Row 10
$ws4.Range("C10").Value = "OrdersGUIController, 54"
$ws4.Range("D10").Value = "Instance methods should not write to ""static"" fields (java:S2696)"
$ws4.Range("E10").Value = "public void setTotalAmount(double totalAmount)"
$ws4.Range("F10").Value = "public static void setTotalAmount(double totalAmount)"
$ws4.Rows.Item(10).RowHeight = 60.1

# Row 11
$ws4.Range("C11").Value = "OrdersGUIController, 59"
$ws4.Range("D11").Value = "Private fields only used as local variables in methods should become local variables (java:S1450)"
$ws4.Range("E11").Value = "    private int tableNumber;"
$ws4.Range("F11").Value = "SonarQube is wrong, variable is specific to instance"
$ws4.Rows.Item(11).RowHeight = 28.8

# Row 12
$ws4.Range("C12").Value = "OrdersGUIController, 63"
$ws4.Range("D12").Value = "Private fields only used as local variables in methods should become local variables (java:S1450)"
$ws4.Range("E12").Value = "private ObservableList<MenuDataModel> menuData;"
$ws4.Range("F12").Value = "SonarQube is wrong, variable is specific to instance"
$ws4.Rows.Item(12).RowHeight = 45.1

# Row 13
$ws4.Range("C13").Value = "OrdersGUIController, 62"
$ws4.Range("D13").Value = "The diamond operator (""<>"") should be used (java:S2293)"
$ws4.Range("E13").Value = "    private TableView<MenuDataModel> table = new TableView<MenuDataModel>();"
$ws4.Range("F13").Value = "    private TableView<MenuDataModel> table = new TableView<>();"
$ws4.Rows.Item(13).RowHeight = 75.15

$ws4.Range("I13").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "Architect. Design Phase Defects" -- new architectural defect row
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Architect. Design Phase Defects")

$ws2.Range("C17").Value = "A08"
$ws2.Range("D17").Value = "n/a"
$ws2.Range("E17").Value = "No class descriptions (comments) are provided"

# This sheet becomes the active tab / active sheet of the workbook, with the
# selection left on F17.
$ws2.Activate()
$ws2.Range("F17").Select()
